# Refresh market-board derived price/profit figures (currentAveragePrice*, LevePrice*, LeveProfit*)
# for a batch of Leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
# Source values were produced by the external scheduled price-scraper run; this script just
# replays the resulting H:N cell updates (some rows gain or lose trailing columns when a price
# leg becomes unavailable / becomes available again).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!43 - Growing Is Knowing / Growth Formula Gamma
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").ClearContents() | Out-Null
$ws.Range("M43").ClearContents() | Out-Null
$ws.Range("N43").Value = 0

# ALC!92 - Whinier than the Sword / Enchanted Koppranickel Ink
$ws.Range("H92").Value = 187.44444
$ws.Range("I92").Value = 235.33333
$ws.Range("J92").Value = 91.666664
$ws.Range("K92").Value = 235.33333
$ws.Range("L92").Value = 91.666664
$ws.Range("M92").Value = 1012.66667
$ws.Range("N92").Value = -2587.666664

# ALC!107 - Another Man's Ink / Enchanted Truegold Ink
$ws.Range("H107").Value = 509.75
$ws.Range("I107").Value = 646.3333
$ws.Range("J107").Value = 100
$ws.Range("K107").Value = 646.3333
$ws.Range("L107").Value = 100
$ws.Range("M107").Value = 1273.6667
$ws.Range("N107").Value = -3940

# ALC!129 - Practical Command / Commanding Craftsman's Draught
$ws.Range("H129").Value = 4624.25
$ws.Range("I129").Value = 1165.6666
$ws.Range("J129").Value = 15000
$ws.Range("K129").Value = 3496.9998
$ws.Range("L129").Value = 45000
$ws.Range("M129").Value = 1503.0002
$ws.Range("N129").Value = -55000

# ALC!135 - For Tired Minds / Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 3935.875
$ws.Range("I135").Value = 3935.875
$ws.Range("K135").Value = 35422.875
$ws.Range("M135").Value = -32887.875

$ws = $wb.Worksheets.Item("ARM")
# ARM!88 - The Mast Chance / Adamantite Rivets
$ws.Range("H88").Value = 600
$ws.Range("I88").Value = 600
$ws.Range("K88").Value = 600
$ws.Range("M88").Value = -194

# ARM!91 - The Rose and the Riveter (L) / Adamantite Rivets
$ws.Range("H91").Value = 600
$ws.Range("I91").Value = 600
$ws.Range("K91").Value = 600
$ws.Range("M91").Value = 804

$ws = $wb.Worksheets.Item("BSM")
# BSM!80 - Unbreaker / Titanium Ingot
$ws.Range("H80").Value = 1045.6
$ws.Range("I80").Value = 1060.5
$ws.Range("K80").Value = 1060.5
$ws.Range("M80").Value = -62.5

# BSM!83 - Attack on Titanium (L) / Titanium Ingot
$ws.Range("H83").Value = 1045.6
$ws.Range("I83").Value = 1060.5
$ws.Range("K83").Value = 5302.5
$ws.Range("M83").Value = -310.5

# BSM!134 - Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 3668.5
$ws.Range("I134").Value = 3668.5
$ws.Range("K134").Value = 11005.5
$ws.Range("M134").Value = -8470.5

$ws = $wb.Worksheets.Item("CRP")
# CRP!26 - As the Worm Turns / Yew Radical
$ws.Range("H26").Value = 2933.3333
$ws.Range("I26").Value = 800
$ws.Range("J26").Value = 4000
$ws.Range("K26").Value = 800
$ws.Range("L26").Value = 4000
$ws.Range("M26").Value = -513
$ws.Range("N26").Value = -4574

# CRP!31 - Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 3567.2222
$ws.Range("I31").Value = 2343.75
$ws.Range("J31").Value = 4546
$ws.Range("K31").Value = 2343.75
$ws.Range("L31").Value = 4546
$ws.Range("M31").Value = -2048.75
$ws.Range("N31").Value = -5136

# CRP!34 - Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 3567.2222
$ws.Range("I34").Value = 2343.75
$ws.Range("J34").Value = 4546
$ws.Range("K34").Value = 2343.75
$ws.Range("L34").Value = 4546
$ws.Range("M34").Value = -2141.75
$ws.Range("N34").Value = -4950

# CRP!58 - You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 3670.3333
$ws.Range("I58").Value = 1012
$ws.Range("J58").Value = 4999.5
$ws.Range("K58").Value = 1012
$ws.Range("L58").Value = 4999.5
$ws.Range("M58").Value = -809
$ws.Range("N58").Value = -5405.5

# CRP!99 - O Pine / Pine Lumber
$ws.Range("H99").Value = 7179.4
$ws.Range("J99").Value = 7199
$ws.Range("L99").Value = 7199
$ws.Range("N99").Value = -10195

# CRP!126 - A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 7179.4
$ws.Range("J126").Value = 7199
$ws.Range("L126").Value = 21597
$ws.Range("N126").Value = -26537

# CRP!132 - Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 6479.8
$ws.Range("I132").Value = 6479.8
$ws.Range("K132").Value = 19439.4
$ws.Range("M132").Value = -16909.4

# CRP!134 - Wood You Be Quiet / Ceiba Lumber
$ws.Range("I134").Value = 1006
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3018
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents() | Out-Null
$ws.Range("N134").Value = -483

# CRP!136 - Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 3670.3333
$ws.Range("I136").Value = 1012
$ws.Range("J136").Value = 4999.5
$ws.Range("K136").Value = 3036
$ws.Range("L136").Value = 14998.5
$ws.Range("M136").Value = -486
$ws.Range("N136").Value = -20098.5

$ws = $wb.Worksheets.Item("CUL")
# CUL!12 - Butter Me Up / Kukuru Butter
$ws.Range("H12").Value = 168.85715
$ws.Range("I12").Value = 32
$ws.Range("J12").Value = 351.33334
$ws.Range("K12").Value = 96
$ws.Range("L12").Value = 1054.00002
$ws.Range("M12").Value = 77
$ws.Range("N12").Value = -1400.00002

# CUL!139 - Najoothie / Wild Banana Blend
$ws.Range("H139").Value = 1819.2858
$ws.Range("I139").Value = 1455.8334
$ws.Range("K139").Value = 4367.5002
$ws.Range("M139").Value = 772.4997999999996

$ws = $wb.Worksheets.Item("GSM")
# GSM!122 - Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 2982.5
$ws.Range("I122").Value = 2982.5
$ws.Range("K122").Value = 8947.5
$ws.Range("M122").Value = -6497.5

# GSM!132 - On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 2221.6667
$ws.Range("J132").Value = 1798.3334
$ws.Range("L132").Value = 5395.0002
$ws.Range("N132").Value = -10455.0002

# GSM!140 - The Right Rod / Ra'Kaznar Rod
$ws.Range("H140").Value = 44997
$ws.Range("I140").Value = 39998
$ws.Range("J140").Value = 49996
$ws.Range("K140").Value = 39998
$ws.Range("L140").Value = 49996
$ws.Range("M140").Value = -34818
$ws.Range("N140").Value = -60356

$ws = $wb.Worksheets.Item("LTW")
# LTW!22 - Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 882.36365
$ws.Range("I22").Value = 672.4286
$ws.Range("J22").Value = 1249.75
$ws.Range("K22").Value = 672.4286
$ws.Range("L22").Value = 1249.75
$ws.Range("M22").Value = -377.4286
$ws.Range("N22").Value = -1839.75

# LTW!27 - Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 882.36365
$ws.Range("I27").Value = 672.4286
$ws.Range("J27").Value = 1249.75
$ws.Range("K27").Value = 672.4286
$ws.Range("L27").Value = 1249.75
$ws.Range("M27").Value = -565.4286
$ws.Range("N27").Value = -1463.75

# LTW!132 - Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 1598.8
$ws.Range("I132").Value = 1598.8
$ws.Range("K132").Value = 4796.4
$ws.Range("M132").Value = -2266.4

# LTW!139 - Giving Gatherers Their Gear / Gomphotherium Doublet of Gathering
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").ClearContents() | Out-Null
$ws.Range("M139").ClearContents() | Out-Null
$ws.Range("N139").Value = 0

$ws = $wb.Worksheets.Item("WVR")
# WVR!132 - Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 1999
$ws.Range("I132").Value = 1999
$ws.Range("K132").Value = 5997
$ws.Range("M132").Value = -3467
